# "Add article; Add artwork" — journal de travail
# Fills in the week-4 (sprint review + réalisation) entries on the
# "Journal" sheet: rows 105-113 get real data, rows 108-117 get the
# formatting of the existing pattern (row 104) pasted onto them so the
# new rows use the same number formats / fills as the rest of the
# table, and the E column gets a shared "=D-C" duration formula spanning
# E108:E117 (row 105-107 reuse the already-existing shared formula).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal")

# ---------------------------------------------------------------------
# 1) Bring the formatting of rows 108-117 in line with the rest of the
#    table (copy the cell styles from row 104, which already has the
#    A:L pattern used throughout the journal). This also fixes up the
#    "spans" bookkeeping Excel keeps for the surrounding 16-row block.
# ---------------------------------------------------------------------
$ws.Range("A104:L104").Copy() | Out-Null
$ws.Range("A108:L117").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) Row heights for the two rows whose content now wraps onto several
#    lines.
# ---------------------------------------------------------------------
$ws.Rows.Item(107).RowHeight = 30
$ws.Rows.Item(109).RowHeight = 135

# ---------------------------------------------------------------------
# 3) Row 105 — Sprint Review day, "Notation des articles" continued.
# ---------------------------------------------------------------------
$ws.Cells.Item(105,1).Value = 44341
$ws.Cells.Item(105,2).Value = 4
$ws.Cells.Item(105,3).Value = 0.33333333333333331
$ws.Cells.Item(105,4).Value = 0.36805555555555558
$ws.Cells.Item(105,6).Value = "Réalisation"
$ws.Cells.Item(105,7).Value = "Notation des articles"

# Row 106 — Sprint Review.
$ws.Cells.Item(106,1).Value = 44341
$ws.Cells.Item(106,2).Value = 4
$ws.Cells.Item(106,3).Value = 0.36805555555555558
$ws.Cells.Item(106,4).Value = 0.3923611111111111
$ws.Cells.Item(106,6).Value = "Communication"
$ws.Cells.Item(106,7).Value = "Sprint Review"

# Row 107 — Corrections following the sprint review.
$ws.Cells.Item(107,1).Value = 44341
$ws.Cells.Item(107,2).Value = 4
$ws.Cells.Item(107,3).Value = 0.3923611111111111
$ws.Cells.Item(107,4).Value = 0.39930555555555558
$ws.Cells.Item(107,6).Value = "Conception"
$ws.Cells.Item(107,7).Value = "Corrections"
$ws.Cells.Item(107,8).Value = "Correstions en fonction de la sprint review: date, historique, MCD, stratégire de test."

# ---------------------------------------------------------------------
# 4) Row 108 — back to "Notation des articles".
# ---------------------------------------------------------------------
$ws.Cells.Item(108,1).Value = 44341
$ws.Cells.Item(108,2).Value = 4
$ws.Cells.Item(108,3).Value = 0.40972222222222227
$ws.Cells.Item(108,4).Value = 0.44097222222222227
$ws.Cells.Item(108,6).Value = "Réalisation"
$ws.Cells.Item(108,7).Value = "Notation des articles"

# Row 109 — Création d'article (documentation links are filled in
# further down, to reproduce the shared-string creation order).
$ws.Cells.Item(109,1).Value = 44341
$ws.Cells.Item(109,2).Value = 4
$ws.Cells.Item(109,3).Value = 0.44097222222222227
$ws.Cells.Item(109,4).Value = 0.4826388888888889
$ws.Cells.Item(109,6).Value = "Réalisation"
$ws.Cells.Item(109,7).Value = "Création d'article"

# Row 110 — Gestion des œuvres / Création d'œuvre.
$ws.Cells.Item(110,1).Value = 44341
$ws.Cells.Item(110,2).Value = 4
$ws.Cells.Item(110,3).Value = 0.4826388888888889
$ws.Cells.Item(110,4).Value = 0.51041666666666663
$ws.Cells.Item(110,6).Value = "Réalisation"
$ws.Cells.Item(110,7).Value = "Gestion des œuvres"
$ws.Cells.Item(110,8).Value = "Création d'œuvre"

# Row 111 — Création d'article continued.
$ws.Cells.Item(111,1).Value = 44341
$ws.Cells.Item(111,2).Value = 4
$ws.Cells.Item(111,3).Value = 0.5625
$ws.Cells.Item(111,4).Value = 0.59375
$ws.Cells.Item(111,6).Value = "Réalisation"
$ws.Cells.Item(111,7).Value = "Création d'article"

# Row 112 — Gestion des œuvres / Création d'œuvre continued.
$ws.Cells.Item(112,1).Value = 44341
$ws.Cells.Item(112,2).Value = 4
$ws.Cells.Item(112,3).Value = 0.59375
$ws.Cells.Item(112,4).Value = 0.6069444444444444
$ws.Cells.Item(112,6).Value = "Réalisation"
$ws.Cells.Item(112,7).Value = "Gestion des œuvres"
$ws.Cells.Item(112,8).Value = "Création d'œuvre"

# Row 113 — Gestion des œuvres / Modification (end of day, no end time).
$ws.Cells.Item(113,1).Value = 44341
$ws.Cells.Item(113,2).Value = 4
$ws.Cells.Item(113,3).Value = 0.6069444444444444
$ws.Cells.Item(113,6).Value = "Réalisation"
$ws.Cells.Item(113,7).Value = "Gestion des œuvres"
$ws.Cells.Item(113,8).Value = "Modification"

# Row 109's documentation links, added last (matches the shared-string
# ordering seen in the target workbook).
$ws.Cells.Item(109,11).Value = "https://www.php.net/manual/fr/reserved.variables.files.php`nhttps://www.php.net/manual/fr/function.strlen.php`nhttps://www.php.net/manual/fr/function.str-contains.php`nhttps://www.php.net/manual/fr/language.operators.comparison.php`n"

# ---------------------------------------------------------------------
# 5) Duration formula: E108:E117 share a single "=D-C" formula, same
#    shape as the existing E103:E107 shared formula.
# ---------------------------------------------------------------------
$ws.Range("E108:E117").Formula = "=D108-C108"

# ---------------------------------------------------------------------
# 6) Scroll position / selection, matching where the author ended up.
# ---------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 97
$ws.Range("G121").Select()
